$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'245.90"
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.Value = "'22.01"
$c.Style = "Normal"

$c = $ws.Range("D4")
$c.Value = "'5.446"
$c.Style = "Normal"

$c = $ws.Range("D5")
$c.Value = "'0.05773"
$c.Style = "Normal"

$c = $ws.Range("D6")
$c.Value = "'3.420"
$c.Style = "Normal"

$c = $ws.Range("D7")
$c.Value = "'6.341"
$c.Style = "Normal"

$c = $ws.Range("D8")
$c.Value = "'0.8187"
$c.Style = "Normal"

$c = $ws.Range("D9")
$c.Value = "'1.036"
$c.Style = "Normal"
$ws.Range("E9").Value = '8FTXTokenFTTBestin24h'

$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$c = $ws.Range("D10")
$c.Value = "'0.1429"
$c.Style = "Normal"
$ws.Range("E10").Value = '9WazirXWRX'

$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$c = $ws.Range("D11")
$c.Value = "'0.07292"
$c.Style = "Normal"
$ws.Range("E11").Value = '10MandalaExchangeTokenMDX'

$ws.Range("B12").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C12").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$c = $ws.Range("D12")
$c.Value = "'0.03123"
$c.Style = "Normal"
$ws.Range("E12").Value = '11LiechtensteinCryptoassetsExchangeLCX'

$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$c = $ws.Range("D13")
$c.Value = "'0.03119"
$c.Style = "Normal"
$ws.Range("E13").Value = '12BitrueCoinBTR'

$ws.Range("B14").Value = 'MCDex'
$ws.Range("C14").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$c = $ws.Range("D14")
$c.Value = "'4.152"
$c.Style = "Normal"
$ws.Range("E14").Value = '13MCDexMCB'

$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$c = $ws.Range("D15")
$c.Value = "'0.09389"
$c.Style = "Normal"
$ws.Range("E15").Value = '14BitMartTokenBMX'

$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$c = $ws.Range("D16")
$c.Value = "'0.001606"
$c.Style = "Normal"
$ws.Range("E16").Value = '15BitForexTokenBF'

$ws.Range("B17").Value = 'CoinExToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$c = $ws.Range("D17")
$c.Value = "'0.04808"
$c.Style = "Normal"
$ws.Range("E17").Value = '16CoinExTokenCET'

$ws.Range("B18").Value = 'One'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$c = $ws.Range("D18")
$c.Value = "'0.0005842"
$c.Style = "Normal"
$ws.Range("E18").Value = '17OneONE'

$c = $ws.Range("D19")
$c.Value = "'0.006252"
$c.Style = "Normal"

$c = $ws.Range("D20")
$c.Value = "'0.004136"
$c.Style = "Normal"

$c = $ws.Range("D21")
$c.Value = "'0.0009919"
$c.Style = "Normal"

$c = $ws.Range("D22")
$c.Value = "'0.0001499"
$c.Style = "Normal"

$c = $ws.Range("D23")
$c.Value = "'3.740"
$c.Style = "Normal"

$c = $ws.Range("D24")
$c.Value = "'2.194"
$c.Style = "Normal"

$c = $ws.Range("D26")
$c.Value = "'0.1328"
$c.Style = "Normal"

$c = $ws.Range("D27")
$c.Value = "'0.0003993"
$c.Style = "Normal"

$c = $ws.Range("D40")
$c.Value = "'0.03883"
$c.Style = "Normal"

$c = $ws.Range("D41")
$c.Value = "'0.006670"
$c.Style = "Normal"

$c = $ws.Range("D44")
$c.Value = "'0.006586"
$c.Style = "Normal"

$c = $ws.Range("D45")
$c.Value = "'0.00005602"
$c.Style = "Normal"

$c = $ws.Range("D47")
$c.Value = "'0.3895"
$c.Style = "Normal"
